$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

function Set-TextValue($range, [string]$text) {
    # Plain assignment is fine for ordinary text, but Excel auto-coerces
    # bare "True"/"False" (and similar look-alike) strings to booleans.
    # Route those through a formula + paste-values round trip so the
    # stored cell stays a literal shared string, not a boolean.
    if ($text -eq "True" -or $text -eq "False") {
        $range.Formula = "=""" + $text + """"
        $range.Copy() | Out-Null
        $range.PasteSpecial(-4163) | Out-Null
        $excel.CutCopyMode = $false
    } else {
        $range.Value = $text
    }
}

# --- Overview sheet: row 3 is the "b.md" file ---
Set-TextValue $overview.Range("E3") "Ready for handoff"
Set-TextValue $overview.Range("F3") "Ready for handoff"
Set-TextValue $overview.Range("G3") "2016-09-06 18:50:38"

# --- zh-cn sheet: row 3 is the "b.md" file ---
Set-TextValue $zhcn.Range("C3") "Ready for handoff"
Set-TextValue $zhcn.Range("F3") "False"
Set-TextValue $zhcn.Range("G3") "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
Set-TextValue $zhcn.Range("H3") "2016-09-06 18:50:33"
Set-TextValue $zhcn.Range("P3") "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/53aac2aaa38f440b4e6ec596113e7f4c6cdf7a31/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c41fc5ad85d27ae85afda87046d2f06f2cce2778/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 40

# --- de-de sheet: row 3 is the "b.md" file ---
Set-TextValue $dede.Range("C3") "Ready for handoff"
Set-TextValue $dede.Range("F3") "False"
Set-TextValue $dede.Range("G3") "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
Set-TextValue $dede.Range("H3") "2016-09-06 18:50:38"
Set-TextValue $dede.Range("P3") "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/53aac2aaa38f440b4e6ec596113e7f4c6cdf7a31/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c41fc5ad85d27ae85afda87046d2f06f2cce2778/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 40
